$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original number format/style so we can force text entry
# for numeric-looking values (avoids Excel auto-converting "543.96" -> 543.96 number)
# across the whole data range, then restore the original style afterwards.
$origStyle = $ws.Range("D2:E51").Style
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.591.11"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.526.82"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "543.96"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "145.84"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "2.551.01"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "5.58"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "2.974.47"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "23.57"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "59.501.57"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "2.534.74"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "326.94"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "62.32"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "0.438"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "0.992"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "6.80"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -8.87%  "
$ws.Range("D33").Value = "1.50"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").Value = "161.23"
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D36").Value = "18.76"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("D39").Value = "37.14"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "5.64"
$ws.Range("E40").Value = "  -7.43%  "
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "297.70"
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "0.993"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.609"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "10.80"
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0939"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "18.95"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "123.69"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("E51").Value = "  -1.79%  "

# Restore original style/number format on the data range
$ws.Range("D2:E51").Style = $origStyle

Write-Output "Applied 87 cell updates"
